# Scheduled-runner refresh of the Adamantoise_Profits leve-crafting
# profit tables: recompute/overwrite cached market-price and profit
# columns (H:N) for the rows whose source prices changed, one sheet
# (crafting class) at a time.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 5833.6665
$ws.Range("J17").Value = 5942.654
$ws.Range("L17").Value = 17827.962
$ws.Range("N17").Value = -18163.962
# Row 28
$ws.Range("H28").Value = 134363.06
$ws.Range("I28").Value = 143781.86
$ws.Range("K28").Value = 143781.86
$ws.Range("M28").Value = -143296.86
# Row 70
$ws.Range("H70").Value = 22732150
$ws.Range("I70").Value = 2331
$ws.Range("K70").Value = 6993
$ws.Range("M70").Value = -6723
# Row 73
$ws.Range("H73").Value = 22732150
$ws.Range("I73").Value = 2331
$ws.Range("K73").Value = 6993
$ws.Range("M73").Value = -6057
# Row 76
$ws.Range("H76").Value = 11494.4
$ws.Range("J76").Value = 9495
$ws.Range("L76").Value = 9495
$ws.Range("N76").Value = -10125
# Row 79
$ws.Range("H79").Value = 11494.4
$ws.Range("J79").Value = 9495
$ws.Range("L79").Value = 9495
$ws.Range("N79").Value = -11679
# Row 116
$ws.Range("H116").Value = 16224.647
$ws.Range("I116").Value = 19355
$ws.Range("K116").Value = 19355
$ws.Range("M116").Value = -15913
# Row 132
$ws.Range("H132").Value = 1425.4531
$ws.Range("I132").Value = 1428.5636
$ws.Range("K132").Value = 4285.6908
$ws.Range("M132").Value = -1755.6908

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 22793084
$ws.Range("I32").Value = 25466472
$ws.Range("K32").Value = 25466472
$ws.Range("M32").Value = -25466185
# Row 110
$ws.Range("H110").Value = 1983.8823
$ws.Range("I110").Value = 1935.4
$ws.Range("J110").Value = 2347.5
$ws.Range("K110").Value = 1935.4
$ws.Range("L110").Value = 2347.5
$ws.Range("M110").Value = 109.5999999999999
$ws.Range("N110").Value = -6437.5
# Row 141
$ws.Range("H141").Value = 649999.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 649999.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 649999.5
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -660359.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4595.6875
$ws.Range("I31").Value = 2256.0833
$ws.Range("K31").Value = 2256.0833
$ws.Range("M31").Value = -1961.0833
# Row 34
$ws.Range("H34").Value = 4595.6875
$ws.Range("I34").Value = 2256.0833
$ws.Range("K34").Value = 2256.0833
$ws.Range("M34").Value = -2054.0833
# Row 58
$ws.Range("H58").Value = 2429.2322
$ws.Range("I58").Value = 2179.6223
$ws.Range("K58").Value = 2179.6223
$ws.Range("M58").Value = -1976.6223
# Row 94
$ws.Range("H94").Value = 3016.75
$ws.Range("J94").Value = 3016.75
$ws.Range("L94").Value = 3016.75
$ws.Range("N94").Value = -3918.75
# Row 136
$ws.Range("H136").Value = 2429.2322
$ws.Range("I136").Value = 2179.6223
$ws.Range("K136").Value = 6538.8669
$ws.Range("M136").Value = -3988.8669

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 296.91666
$ws.Range("J7").Value = 585.25
$ws.Range("L7").Value = 1755.75
$ws.Range("N7").Value = -1979.75
# Row 12
$ws.Range("H12").Value = 135.125
$ws.Range("J12").Value = 164.33333
$ws.Range("L12").Value = 492.99999
$ws.Range("N12").Value = -838.99999
# Row 121
$ws.Range("H121").Value = 1154.0526
$ws.Range("J121").Value = 1318.2307
$ws.Range("L121").Value = 3954.6921
$ws.Range("N121").Value = -6574.6921
# Row 131
$ws.Range("H131").Value = 1788.5333
$ws.Range("J131").Value = 1789.8148
$ws.Range("L131").Value = 5369.4444
$ws.Range("N131").Value = -15449.4444
# Row 132
$ws.Range("H132").Value = 591681.9
$ws.Range("I132").Value = 2755.4443
$ws.Range("J132").Value = 1254224.1
$ws.Range("K132").Value = 24798.9987
$ws.Range("L132").Value = 11288016.9
$ws.Range("M132").Value = -22268.9987
$ws.Range("N132").Value = -11293076.9

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 14354.357
$ws.Range("I2").Value = 80.083336
$ws.Range("K2").Value = 80.083336
$ws.Range("M2").Value = 32.916664
# Row 7
$ws.Range("H7").Value = 5256625
$ws.Range("I7").Value = 5259500.5
$ws.Range("J7").Value = 5253750
$ws.Range("K7").Value = 5259500.5
$ws.Range("L7").Value = 5253750
$ws.Range("M7").Value = -5259388.5
$ws.Range("N7").Value = -5253974
# Row 8
$ws.Range("H8").Value = 5256625
$ws.Range("I8").Value = 5259500.5
$ws.Range("J8").Value = 5253750
$ws.Range("K8").Value = 5259500.5
$ws.Range("L8").Value = 5253750
$ws.Range("M8").Value = -5259361.5
$ws.Range("N8").Value = -5254028
# Row 44
$ws.Range("H44").Value = 99990
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
# Row 54
$ws.Range("H54").Value = 59450
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 59450
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 59450
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -60230
# Row 132
$ws.Range("H132").Value = 1613.4546
$ws.Range("I132").Value = 789.1667
$ws.Range("K132").Value = 2367.5001
$ws.Range("M132").Value = 162.4998999999998

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1071.5714
$ws.Range("J22").Value = 1460.5
$ws.Range("L22").Value = 1460.5
$ws.Range("N22").Value = -2050.5
# Row 23
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
# Row 24
$ws.Range("H24").Value = 13242.667
$ws.Range("I24").Value = 14291.2
$ws.Range("J24").Value = 8000
$ws.Range("K24").Value = 14291.2
$ws.Range("L24").Value = 8000
$ws.Range("M24").Value = -13948.2
$ws.Range("N24").Value = -8686
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
# Row 27
$ws.Range("H27").Value = 1071.5714
$ws.Range("J27").Value = 1460.5
$ws.Range("L27").Value = 1460.5
$ws.Range("N27").Value = -1674.5
# Row 45
$ws.Range("H45").Value = 29540.5
$ws.Range("I45").Value = 29540.5
$ws.Range("K45").Value = 29540.5
$ws.Range("M45").Value = -29133.5
# Row 48
$ws.Range("H48").Value = 25707.334
$ws.Range("I48").Value = 25707.334
$ws.Range("K48").Value = 25707.334
$ws.Range("M48").Value = -25046.334
# Row 82
$ws.Range("H82").Value = 2465.8276
$ws.Range("J82").Value = 2513.2222
$ws.Range("L82").Value = 2513.2222
$ws.Range("N82").Value = -3235.2222
# Row 85
$ws.Range("H85").Value = 2465.8276
$ws.Range("J85").Value = 2513.2222
$ws.Range("L85").Value = 2513.2222
$ws.Range("N85").Value = -5009.2222

$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 9544.817999999999
$ws.Range("I2").Value = 9499.4
$ws.Range("J2").Value = 9999
$ws.Range("K2").Value = 9499.4
$ws.Range("L2").Value = 9999
$ws.Range("M2").Value = -9387.4
$ws.Range("N2").Value = -10223
# Row 4
$ws.Range("H4").Value = 23365.455
$ws.Range("I4").Value = 25666
$ws.Range("J4").Value = 360
$ws.Range("K4").Value = 25666
$ws.Range("L4").Value = 360
$ws.Range("M4").Value = -25553
$ws.Range("N4").Value = -586
# Row 74
$ws.Range("H74").Value = 58540.668
$ws.Range("J74").Value = 58540.668
$ws.Range("L74").Value = 58540.668
$ws.Range("N74").Value = -60412.668
# Row 77
$ws.Range("H77").Value = 58540.668
$ws.Range("J77").Value = 58540.668
$ws.Range("L77").Value = 175622.004
$ws.Range("N77").Value = -184982.004
# Row 105
$ws.Range("H105").Value = 37807.5
$ws.Range("J105").Value = 37807.5
$ws.Range("L105").Value = 37807.5
$ws.Range("N105").Value = -44795.5
